# Updates LR-pair metrics (columns G:T) for each row with newly computed TPM-based
# statistics, and re-labels the Target cluster (column D) for several rows from
# "Resolving-Mac" to "MuSCs" to reflect the corrected cluster assignment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Target cluster label for the affected rows
foreach ($r in @(4,7,10,13,16,19)) {
    $ws.Range("D$r").Value = "MuSCs"
}

# New G:T metrics per row (Ligand/Receptor expression + specificity + edge weights)
$rowData = @{}
$rowData[2] = @(58.45400900000001,175.362027,0.654344352429287,0.6543443524292871,2,0.6666666666666666,0.8595406666666667,2.578622,0.7048880036936188,0.7048880036936188,50.24359786519934,452.192380786794,0.4612394843120739,0.461239484312074)
$rowData[3] = @(58.45400900000001,175.362027,0.654344352429287,0.6543443524292871,2,0.6666666666666666,0.282763,0.8482890000000001,0.2318869302151522,0.2318869302151522,16.528630946867,148.757678521803,0.151733903188449,0.1517339031884491)
$rowData[4] = @(58.45400900000001,175.362027,0.654344352429287,0.6543443524292871,1,0.3333333333333333,0.07709666666666666,0.23129,0.06322506609122898,0.06322506609122898,4.506609247203333,40.55948322483,0.0413709649287641,0.0413709649287641)
$rowData[5] = @(18.46163566666667,55.384907,0.2066627634571724,0.2066627634571724,2,0.6666666666666666,0.8595406666666667,2.578622,0.7048880036936188,0.7048880036936188,15.86852662868378,142.816739658154,0.1456741027711328,0.1456741027711328)
$rowData[6] = @(18.46163566666667,55.384907,0.2066627634571724,0.2066627634571724,2,0.6666666666666666,0.282763,0.8482890000000001,0.2318869302151522,0.2318869302151522,5.220267486013667,46.982407374123,0.04792239380786385,0.04792239380786385)
$rowData[7] = @(18.46163566666667,55.384907,0.2066627634571724,0.2066627634571724,1,0.3333333333333333,0.07709666666666666,0.23129,0.06322506609122898,0.06322506609122898,1.423330571114444,12.80997514003,0.01306626687817575,0.01306626687817575)
$rowData[8] = @(10.312204,30.936612,0.1154366067261307,0.1154366067261307,2,0.6666666666666666,0.8595406666666667,2.578622,0.7048880036936188,0.7048880036936188,8.863758700962666,79.77382830866401,0.08136987926834761,0.08136987926834763)
$rowData[9] = @(10.312204,30.936612,0.1154366067261307,0.1154366067261307,2,0.6666666666666666,0.282763,0.8482890000000001,0.2318869302151522,0.2318869302151522,2.915909739652,26.243187656868,0.02676824036817623,0.02676824036817623)
$rowData[10] = @(10.312204,30.936612,0.1154366067261307,0.1154366067261307,1,0.3333333333333333,0.07709666666666666,0.23129,0.06322506609122898,0.06322506609122898,0.7950365543866665,7.15532898948,0.007298487089606819,0.00729848708960682)
$rowData[11] = @(0.2278333333333333,0.6835,0.002550405994596639,0.002550405994596639,2,0.6666666666666666,0.8595406666666667,2.578622,0.7048880036936188,0.7048880036936188,0.1958320152222222,1.762488137,0.001797750590139463,0.001797750590139463)
$rowData[12] = @(0.2278333333333333,0.6835,0.002550405994596639,0.002550405994596639,2,0.6666666666666666,0.282763,0.8482890000000001,0.2318869302151522,0.2318869302151522,0.06442283683333334,0.5798055315,0.0005914058168893367,0.0005914058168893366)
$rowData[13] = @(0.2278333333333333,0.6835,0.002550405994596639,0.002550405994596639,1,0.3333333333333333,0.07709666666666666,0.23129,0.06322506609122898,0.06322506609122898,0.01756519055555555,0.158086715,0.0001612495875678391,0.0001612495875678391)
$rowData[14] = @(1.290132,3.870396,0.01444196219438603,0.01444196219438603,2,0.6666666666666666,0.8595406666666667,2.578622,0.7048880036936188,0.7048880036936188,1.108920919368,9.980288274312,0.01017996590061948,0.01017996590061948)
$rowData[15] = @(1.290132,3.870396,0.01444196219438603,0.01444196219438603,2,0.6666666666666666,0.282763,0.8482890000000001,0.2318869302151522,0.2318869302151522,0.364801594716,3.283214352444,0.003348902279539461,0.00334890227953946)
$rowData[16] = @(1.290132,3.870396,0.01444196219438603,0.01444196219438603,1,0.3333333333333333,0.07709666666666666,0.23129,0.06322506609122898,0.06322506609122898,0.09946487676,0.8951838908399999,0.0009130940142270872,0.0009130940142270872)
$rowData[17] = @(0.5863683333333334,1.759105,0.006563909198427097,0.006563909198427097,2,0.6666666666666666,0.8595406666666667,2.578622,0.7048880036936188,0.7048880036936188,0.5040074281455557,4.53606685331,0.004626820851305459,0.004626820851305459)
$rowData[18] = @(0.5863683333333334,1.759105,0.006563909198427097,0.006563909198427097,2,0.6666666666666666,0.282763,0.8482890000000001,0.2318869302151522,0.2318869302151522,0.1658032690383334,1.492229421345,0.00152208475423426,0.00152208475423426)
$rowData[19] = @(0.5863683333333334,1.759105,0.006563909198427097,0.006563909198427097,1,0.3333333333333333,0.07709666666666666,0.23129,0.06322506609122898,0.06322506609122898,0.04520704393888889,0.40686339545,0.0004150035928873791,0.0004150035928873791)

foreach ($r in $rowData.Keys) {
    $arr = New-Object "object[,]" 1,14
    for ($i = 0; $i -lt 14; $i++) { $arr[0,$i] = $rowData[$r][$i] }
    $ws.Range("G$r`:T$r").Value = $arr
}
